# Applies the updated numeric values (columns E:AC, rows 2:29) from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$firstCol = 5  # column E
$numRows = 28
$numCols = 25  # E..AC

# One row per spreadsheet row (r=2..29); each inner array holds the 25 new values for columns E..AC, as strings so scientific notation parses exactly as a double.
$newValues = @(
    @("0.1588460368111074", "0.07700054126458425", "0.04972798644600025", "0.07891296682397321", "0.004380961057758153", "0.03870225026195415", "0.01962776015267511", "0.08236898882062581", "0.0008228562329691285", "0.006282670095193112", "0.01257273708299272", "0.06328685618531611", "0.06095752621508826", "0.007458641499828065", "0.04423529301198196", "0.0326590288264908", "0.002154843682456897", "0.01585982036200003", "0.1008112039719405", "0.05341467403828857", "0.00952643655581298", "0.1002752540763608", "0.08806019550859796", "0.0509005078271112", "0.1569408081499378"),
    @("0.2146249116276946", "0.005768524047347839", "0.05448111263446578", "0.01442016154170439", "0.02902059526282754", "0.02250335032961909", "0.1125850842054206", "0.1048080162660824", "0.001853811274807991", "0.003993392845164164", "0.1136984414110694", "0.004139404909854346", "0.1170203456825184", "0.002062994291711913", "0.06952624562355333", "0.005178957131098642", "0.01001742258942492", "0.004483639728113275", "0.05612089987547601", "0.05875595914069325", "0.007933238696644135", "0.06116855078521913", "0.02850989158918382", "0.1119499601379995", "0.1717004433804159"),
    @("0.2414158924066542", "0.06709450760274706", "0.0330218666434625", "0.05794682377151639", "0.08243176751657431", "0.02796525219568386", "0.05482227884127638", "0.04416213677064301", "0.09018795740462954", "0.009304941303576267", "0.05404527588718264", "0.0004428004940759327", "0.003560887472620206", "0.0875174807645262", "0.06571211614052604", "0.04025677613576578", "0.08896372502260713", "0.03613959295713459", "0.01228409621584211", "0.03924645219450483", "0.05735083706540997", "0.03741092734539427", "0.002798171620287707", "0.00733332863401339", "-0.6422546796980082"),
    @("0.4216481654301224", "0.04995640328017423", "0.000869889394924528", "0.1050809500553407", "0.06242042389665681", "0.01266172361660683", "0.0125073283454852", "0.04185771435499722", "0.04595176236426114", "0.009325478957632211", "0.08316155265763188", "0.02503956365889011", "0.01639285340469301", "0.06665974364163922", "0.08972098872025652", "0.02192478768317018", "0.0267470785683654", "0.01769669743866899", "0.09565174128373113", "0.08824341634299022", "0.001089740923809289", "0.03204073317808369", "0.05279107648510264", "0.04220835174688892", "-0.5589303733638036"),
    @("0.2445577055276964", "0.01783264882109302", "0.03660700401295174", "0.04714634984098481", "0.04719303260530715", "0.01824818723007508", "0.04726248123106434", "0.02956850194039538", "0.006613766441415303", "0.00750643011018602", "0.0705536135084701", "0.02919688161923314", "0.001522249872934383", "0.09265920022748386", "0.08372404314027809", "0.0841202936768333", "0.06307699824116532", "0.006364783661002463", "0.06497695041795105", "0.05979230802497879", "0.07261650096141538", "0.08837948564438874", "2.49680374780602e-05", "0.0250133207329144", "0.2121750077080789"),
    @("0.3992962682419841", "0.04153765137902524", "0.07547310951908325", "0.01632008483007047", "0.01070826771858457", "0.04847140383878051", "0.0006145732598764111", "0.008039112201915074", "0.01925116995354929", "0.08981755891302001", "0.01688734509498854", "0.02507617361277198", "0.009320457998957012", "0.06610205342160813", "0.02460634933678289", "0.1023247059416084", "0.03086741327397521", "0.1157647316069643", "0.0408386877159185", "0.11011951515794", "0.01820628785282885", "0.0293239078588314", "0.02255388708373476", "0.07777555242918513", "0.3841434140024276"),
    @("0.2740557721257694", "0.08379440209804201", "0.06349486489963745", "0.003975972204747372", "0.05321830558491093", "0.02293646181787647", "0.03219817993159932", "0.06539306895552809", "0.02126744649506489", "0.008950189011616353", "0.01179297204233967", "0.04156094237426075", "0.06891367367895589", "0.08249266972201216", "0.066621330436459", "0.04143480102532381", "0.0005747150652272894", "0.01444966519773094", "0.07025317187503471", "0.02104311005062241", "0.06723129184927314", "0.0601599462356228", "0.05588008247567046", "0.04236273697244423", "0.1118222900081235"),
    @("0.3355899301812887", "0.04852206181110189", "0.02542373501344083", "0.01583863242988653", "0.02344387327218946", "0.09878917902809252", "0.03141586771655899", "0.01360524532508345", "0.005717079108322028", "0.07874105726292763", "0.04311953348016985", "0.07045585854531446", "0.009876593555030977", "0.08525191941163818", "0.08079063712573198", "0.07456675238972818", "0.0004459577790353026", "0.09911697171726323", "0.0116092253402012", "0.05573538425739202", "0.01568752099988377", "0.004237541191161499", "0.02909937542927203", "0.07850999781057409", "0.2135692896278156"),
    @("0.4238758601859445", "0.04726582972630274", "0.02701142012555267", "0.02066285962297705", "0.05497966620109976", "0.09896398350601619", "0.01373978230215595", "0.01487473841518705", "0.1214243649114135", "0.0177807634994685", "0.004038936997248371", "0.05888084344021671", "0.07607566625839321", "0.01046971816030925", "0.01069961023313779", "0.0944280394224243", "0.01292986892127918", "0.02274785193972829", "0.1074428290570988", "0.03120586977519829", "0.1244500368766979", "0.0002675449636102513", "0.003368414356670089", "0.02629136128781386", "-0.1948932169126427"),
    @("1.177415973328571", "0.07424310272793004", "0.005622174115404762", "0.09579603580551516", "0.02755119891955621", "0.05022860441824364", "0.02647135829326475", "0.08718006133063158", "0.01695648418708114", "0.07942499835431642", "0.08209463844585704", "0.04531197878801218", "0.04189004420591425", "0.04480420440678532", "0.07093116681179201", "0.02520014076483553", "0.007254742027225473", "0.05504216162204528", "0.01985495775733922", "0.05890523247605881", "0.02035656335912677", "0.02472428537294491", "0.02468434520333521", "0.01547152060678427", "-0.1431703458337011"),
    @("0.1848448830903597", "0.05173426071429667", "0.02776646555623737", "0.038707172035578", "0.04537481763721113", "0.01983141006761854", "0.007932308128964571", "0.02595057643955742", "0.08442783284107352", "0.03648531882698258", "0.08335723234905769", "0.08759491969023188", "0.01499005452169955", "0.07703490147842193", "0.03144703728561646", "0.03696926880159203", "0.04511064277960807", "0.0754346985005634", "0.06087291525431412", "0.02294817559667535", "0.006514808316185531", "0.08308946633161104", "0.03015672013062892", "0.006268996716274089", "0.1775940895285582"),
    @("0.2384357213542734", "0.1021092640762821", "0.01568049344536434", "0.02723179686018686", "0.003982995542286656", "0.06758954990424852", "0.02197674354171031", "0.006897198112519034", "0.02866474937450238", "0.1152995365122887", "0.1195485554586767", "0.02755715054153832", "0.02378740756304568", "0.1417545326681444", "0.001661896556798629", "0.03260322210148865", "0.002465625509330833", "0.03413080661296766", "0.124082055613827", "0.01194137444604743", "0.01562198152335668", "0.0295034067739726", "0.02808731241034504", "0.01782234485107146", "0.1818157329640047"),
    @("0.2172607176495693", "0.06715955956115509", "0.0203907309063375", "0.03092977683222757", "0.05172797820140758", "0.0854467753722435", "0.001913829968874008", "0.03434946305115567", "0.02274580315773428", "0.07029195272897498", "0.006532480723971011", "0.0707810278848038", "0.08799666646779784", "0.08246715096279905", "0.07851825620456787", "0.02865810585697029", "0.05827579442689038", "0.09019321087285181", "0.0008474876938056129", "0.001789052164160249", "0.02008902454790404", "0.06678168109881515", "0.00661464120988409", "0.01549955010466864", "0.1557858288607722"),
    @("0.3853914660150309", "0.03173680857224954", "0.07996745831300556", "0.02879273227342352", "0.07260047016683903", "0.02401899583757255", "0.09287593220645922", "0.0902030693335206", "0.01650575376790574", "0.002697916000611706", "0.05599393088578417", "0.03737428672260822", "0.01657704236194494", "0.01886302734712945", "0.05220269757796922", "0.09739605488160574", "0.05318164362382365", "0.02771522551105728", "0.003788769602300852", "0.01985226708941638", "0.004412486173148281", "0.08640431069059377", "0.01398577146117246", "0.07285334959985798", "0.2485130745104423"),
    @("0.2986159941942825", "0.08694966650398711", "0.08868518870475982", "0.005620119701126668", "0.00503789066209729", "0.05349542568575064", "0.04859958750089569", "0.04858713144394591", "0.09163859631350343", "0.0576360656865446", "0.09342057014163994", "0.09630066783835661", "0.06054692582949145", "0.04047063338534485", "0.02250286717110209", "0.009875811587936286", "0.002378785702262271", "0.006553083300440846", "0.001274487219490967", "0.02649840920756449", "0.07395726951642309", "0.02436747431855013", "0.02314101348349566", "0.03246232909528991", "-0.007929441062452365"),
    @("0.534122564628467", "0.0631549928064888", "0.03562162990918559", "0.008007969282927411", "0.07139010675618124", "0.03619269875228145", "0.06778473569420056", "0.06913542876179599", "0.0129943550388068", "0.07206073204254218", "0.05716730373378795", "0.07457489726484161", "0.02819778179590006", "0.03482409564037213", "0.007416930347590786", "0.03625700944311842", "0.02993240158520439", "0.01623207512951162", "0.06903542474836248", "0.0444211446962323", "0.007406890298984376", "0.07464024595277224", "0.0744417272295762", "0.009109423089335248", "0.01336789697259665"),
    @("0.1644155421449354", "0.0961185764494758", "0.09968134021584801", "0.01496650052440109", "0.0005600055907574322", "0.04965901197605599", "0.09056931462148687", "0.06806606128401621", "0.02169557971776746", "0.07617166435012898", "0.01230127496698816", "0.07875122455244608", "0.01905050820927658", "0.01543280203578008", "0.002018870027597057", "0.07249014455320162", "0.07667607000811559", "0.03981706090656179", "0.003680966909345237", "0.008421428780834636", "0.0471649129393378", "0.05720915525809191", "0.00624590994297463", "0.0432516161795111", "0.1243236406798293"),
    @("0.7149269615976624", "0.01636667469940909", "0.08814497814339119", "0.08891396576259857", "0.0895686906194344", "0.02747347188078936", "0.03249106576326908", "0.07969694244944063", "0.03243923530143528", "0.005995640610422877", "0.01633637230980549", "0.04727864396048693", "0.0119729912723074", "0.02055014939283486", "0.09174711750208725", "0.08886041396279831", "0.04800691210636884", "0.08159130188909694", "0.01172619759273884", "0.02316702856775708", "0.02217355376929672", "0.05506778668181106", "0.01041072822573241", "0.01002013753668747", "0.1688672931566585"),
    @("0.1631379394963398", "0.08595341572668373", "0.01950064148186779", "0.003474053306131192", "0.05439082231762425", "0.08349825427513287", "0.004635730532881596", "0.006634378374589362", "0.06937705244973701", "0.02463684741669372", "0.03545344960275154", "0.008407701990166418", "0.008890401007688894", "0.07859507994261165", "0.08093761842388969", "0.05453590720190461", "0.06742087517016669", "0.001810051228873713", "0.04248419273841035", "0.0741046198688745", "0.001387135871889993", "0.04220726527106833", "0.08596658129805844", "0.06569792450230331", "0.1195450499911382"),
    @("0.2033909211027964", "0.09735631011492568", "0.01864418349922766", "0.07275364906294811", "0.02565537067790677", "0.05743900417001055", "0.07676104405567288", "0.02564443169201668", "0.09431573568117603", "0.01007019103551254", "0.0294509794724934", "0.01102266802195152", "0.04831583757712255", "0.04794830991483325", "0.09235572789238806", "0.03806463335820089", "0.06294696925498522", "0.05287818307490891", "0.06106760933123404", "0.01140847986338402", "0.02926729071984658", "0.01739191187867974", "0.009131792824123814", "0.0101096868264511", "0.1141721250941084"),
    @("0.2066819518762127", "0.06267575227013374", "0.0703337224653906", "0.03714249137777401", "0.006785894844895174", "0.02353652478778991", "0.06695129105478985", "0.1101242578885858", "0.0818892368099481", "0.04610021339112453", "0.0117234072385253", "0.1159850260805384", "0.02357405821708502", "0.008127346697756475", "0.004406301807919772", "0.007302376012196759", "0.1223854496475494", "0.00938582732981798", "0.08539826306869622", "0.03846611262785683", "0.01124800097420168", "0.05355379472855475", "0.002428670421998168", "0.0004759802568715448", "0.04705371038395056"),
    @("0.8951378135425161", "0.142389275730587", "0.0005443352201751613", "0.05271723435096708", "0.04200492739567125", "0.01849007413829759", "0.02051273734227063", "0.003076665707925182", "0.04914698213105217", "0.00717936157389684", "0.004162922264542885", "0.03035030319914847", "0.0486565395893489", "0.1322984955278149", "0.1342671751139161", "0.01066922082809889", "0.06257443368106246", "0.08780078907076236", "0.02062353107117692", "0.06238478870599565", "0.04035076340987099", "0.01671045272516841", "0.002483871707606279", "0.01060511951464407", "0.07499123292960645"),
    @("0.1600192863746284", "0.05799235321414033", "0.07072621867351971", "0.07546314262189985", "0.03831492679929525", "0.03558824390555457", "0.04856604805858195", "0.0652740001992274", "0.004952621503237649", "0.05453361291076365", "0.002585075111935606", "0.06288460886894147", "0.05563538383100373", "0.07634872720521227", "0.004015512051943155", "0.02697082380528248", "0.003545670894762767", "0.06340504430971224", "0.01817425500536735", "0.04283605745251568", "0.01987772883918981", "0.07477452840879838", "0.04385113927255289", "0.05368427705656171", "-0.2514461386044626"),
    @("0.4624028108164472", "0.1004140339635165", "0.1080640365295426", "0.01849092337343243", "0.004054459631026469", "0.03775954502938665", "0.01920141069312358", "0.03642073455929215", "0.1048347128737377", "0.03817154153188491", "0.004617798279817552", "0.02127467999102571", "0.05055911979961197", "0.06977372166979201", "0.005756841035972762", "0.007243126107799584", "0.02048798705498536", "0.0990020743767136", "0.007247425448053966", "0.06225999419642728", "0.0291163343512485", "0.03950908928029623", "0.0138127383370324", "0.1019276718862801", "-0.1572611063694041"),
    @("0.1188758661652794", "0.07520627666273239", "0.092085065950182", "0.04281728578651923", "0.05196750561932612", "0.004178135019150286", "0.007262030182281898", "0.05201721854792293", "0.01454080082822934", "0.01714778008671847", "0.01191591191985152", "0.09515105742783075", "0.09681481021448747", "0.06156131308025429", "0.004411122000577991", "0.07343581632356586", "0.07708650071714231", "0.0003955241094539645", "0.05140170741569321", "0.1091005206258127", "0.01825152647331393", "0.02183870582317038", "0.001519031937051574", "0.01989435324873151", "0.1693744719191226"),
    @("0.1790994661794387", "0.09850417863244602", "0.03872273149002235", "0.003496313879980348", "0.02902582799052891", "0.006310867667427428", "0.01017443557567802", "0.008026479400368585", "0.1115786319948792", "0.008349908103346342", "0.06297042634229107", "0.1001003365491223", "0.002565527391885512", "0.05037998909810337", "0.0003863539514109807", "0.03649251303127041", "0.1161524708860902", "0.07763669877033115", "0.02109125566571944", "0.02035298803344351", "0.05001700125358111", "0.01468047787060063", "0.0380811990839607", "0.09490338733751241", "0.1959373634708947"),
    @("0.1391776576932503", "0.008756667240066037", "0.01617482849121875", "0.02929000235833514", "0.03097986565520613", "0.06783872041513202", "0.03327036556862934", "0.07164481238039372", "0.06892683364075963", "0.02537566714290543", "0.04181157774922842", "0.06253453515633511", "0.08097271511706775", "0.07488496094910849", "0.07896737384054996", "0.008971082834663342", "0.0497467412099328", "0.05291815156392621", "0.07748292640471725", "0.0185643953695053", "0.02890380403540935", "0.008860070102010717", "0.0179186740720408", "0.04520522870285849", "0.02272725869068091"),
    @("0.1503410094277718", "0.02166425066157233", "0.02376507143596479", "0.0627477525929311", "0.09137703533148168", "0.08181280951052061", "0.003232624141010214", "0.00124347229985029", "0.05929069557287127", "0.06992689409561213", "0.05540022026065503", "0.04185257466020784", "0.02493140962754473", "0.08407476560369083", "0.06399491410708348", "0.002431558861958668", "0.01016870047175106", "0.0256072415480945", "0.08434992177827468", "0.006200425623704009", "0.06591836304943419", "0.07818036672712245", "0.04002624312913208", "0.001802688909532164", "0.0004393102680317901")
)

$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = [double]$newValues[$i][$j]
    }
}

$startCell = $ws.Cells.Item($firstRow, $firstCol)
$endCell = $ws.Cells.Item($firstRow + $numRows - 1, $firstCol + $numCols - 1)
$ws.Range($startCell, $endCell).Value = $data

